# Translation-material classification table: renumber the category
# markers on slide 2 ("翻訳対象の分類（案）：前頁①～③で分類").
#
# - "①原作者の独自性が低い" (top-right header cell) becomes "②原作者の独自性が低い"
# - both "翻訳量" sub-header cells become "②翻訳量"
# - "更新頻度" row header becomes "③更新頻度"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

# Row 1, Col 5: "①原作者の独自性が低い（客観的、技術的説明など）"
# Only the leading "①原作者の独自性が低い" run needs to turn into "②原作者の独自性が低い";
# the remainder of the cell ("（客観的、技術的説明など）") is untouched.
$cell1 = $tbl.Cell(1, 5)
$tr1 = $cell1.Shape.TextFrame.TextRange
$head1 = $tr1.Characters(1, 11)
$head1.Text = "②原作者の独自性が低い"

# Row 2, Col 3: "翻訳量" -> "②翻訳量"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "②翻訳量"

# Row 2, Col 5: "翻訳量" -> "②翻訳量"
$tbl.Cell(2, 5).Shape.TextFrame.TextRange.Text = "②翻訳量"

# Row 4, Col 1: "更新頻度" -> "③更新頻度"
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "③更新頻度"
